$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inputs")

$ws.Range("C16:L16").Value = 7500
